$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.82"
$ws.Range("E2").Value = "'5.45%"
$ws.Range("G2").Value = "'20"
$ws.Range("D3").Value = "'32.04"
$ws.Range("E3").Value = "'9.83%"
$ws.Range("G3").Value = "'20"
$ws.Range("D4").Value = "'5.271"
$ws.Range("E4").Value = "'1.35%"
$ws.Range("G4").Value = "'20"
$ws.Range("D5").Value = "'0.07485"
$ws.Range("E5").Value = "'7.65%"
$ws.Range("G5").Value = "'20"
$ws.Range("D6").Value = "'7.865"
$ws.Range("E6").Value = "'5.61%"
$ws.Range("G6").Value = "'20"
$ws.Range("D7").Value = "'3.819"
$ws.Range("E7").Value = "'7.54%"
$ws.Range("G7").Value = "'20"
$ws.Range("D8").Value = "'1.513"
$ws.Range("E8").Value = "'7.81%"
$ws.Range("G8").Value = "'20"
$ws.Range("D9").Value = "'0.9216"
$ws.Range("E9").Value = "'1.92%"
$ws.Range("G9").Value = "'20"
$ws.Range("D10").Value = "'0.1689"
$ws.Range("E10").Value = "'5.42%"
$ws.Range("G10").Value = "'20"
$ws.Range("D11").Value = "'0.07951"
$ws.Range("E11").Value = "'2.76%"
$ws.Range("G11").Value = "'20"
$ws.Range("D12").Value = "'0.08030"
$ws.Range("E12").Value = "'3.96%"
$ws.Range("G12").Value = "'20"
$ws.Range("D13").Value = "'0.03032"
$ws.Range("E13").Value = "'3.48%"
$ws.Range("G13").Value = "'20"
$ws.Range("D14").Value = "'0.09904"
$ws.Range("E14").Value = "'9.77%"
$ws.Range("G14").Value = "'20"
$ws.Range("D15").Value = "'0.001493"
$ws.Range("E15").Value = "'-6.90%"
$ws.Range("G15").Value = "'20"
$ws.Range("D16").Value = "'0.04600"
$ws.Range("E16").Value = "'1.74%"
$ws.Range("G16").Value = "'20"
$ws.Range("D17").Value = "'0.006442"
$ws.Range("E17").Value = "'3.10%"
$ws.Range("G17").Value = "'20"
$ws.Range("E18").Value = "'-0.29%"
$ws.Range("G18").Value = "'20"
$ws.Range("D19").Value = "'2.232"
$ws.Range("E19").Value = "'0.09%"
$ws.Range("G19").Value = "'20"
$ws.Range("D20").Value = "'0.3302"
$ws.Range("E20").Value = "'2.09%"
$ws.Range("G20").Value = "'20"
$ws.Range("D21").Value = "'0.1345"
$ws.Range("E21").Value = "'0.30%"
$ws.Range("G21").Value = "'20"
$ws.Range("D22").Value = "'4.487"
$ws.Range("E22").Value = "'11.98%"
$ws.Range("G22").Value = "'20"
$ws.Range("D23").Value = "'0.1623"
$ws.Range("E23").Value = "'1.52%"
$ws.Range("G23").Value = "'20"
$ws.Range("E24").Value = "'0.47%"
$ws.Range("G24").Value = "'20"
$ws.Range("D25").Value = "'0.004436"
$ws.Range("E25").Value = "'6.93%"
$ws.Range("G25").Value = "'20"
$ws.Range("E26").Value = "'19.66%"
$ws.Range("G26").Value = "'20"
$ws.Range("D27").Value = "'0.0001778"
$ws.Range("E27").Value = "'6.65%"
$ws.Range("G27").Value = "'20"
$ws.Range("G28").Value = "'20"
$ws.Range("G29").Value = "'20"
$ws.Range("G30").Value = "'20"
$ws.Range("G31").Value = "'20"
$ws.Range("G32").Value = "'20"
$ws.Range("G33").Value = "'20"
$ws.Range("G34").Value = "'20"
$ws.Range("G35").Value = "'20"
$ws.Range("G36").Value = "'20"
$ws.Range("G37").Value = "'20"
$ws.Range("G38").Value = "'20"
$ws.Range("D39").Value = "'0.01697"
$ws.Range("E39").Value = "'2,502.36%"
$ws.Range("G39").Value = "'20"
$ws.Range("D40").Value = "'0.04498"
$ws.Range("E40").Value = "'2.39%"
$ws.Range("G40").Value = "'20"
$ws.Range("D41").Value = "'0.007123"
$ws.Range("E41").Value = "'2.70%"
$ws.Range("G41").Value = "'20"
$ws.Range("D42").Value = "'0.1349"
$ws.Range("E42").Value = "'8.14%"
$ws.Range("G42").Value = "'20"
$ws.Range("D43").Value = "'0.002226"
$ws.Range("E43").Value = "'7.72%"
$ws.Range("G43").Value = "'20"
$ws.Range("D44").Value = "'0.01287"
$ws.Range("E44").Value = "'10.39%"
$ws.Range("G44").Value = "'20"
$ws.Range("D45").Value = "'0.00006183"
$ws.Range("E45").Value = "'6.60%"
$ws.Range("G45").Value = "'20"
$ws.Range("D46").Value = "'0.7093"
$ws.Range("E46").Value = "'-63.23%"
$ws.Range("G46").Value = "'20"
$ws.Range("D47").Value = "'0.01608"
$ws.Range("E47").Value = "'23.71%"
$ws.Range("G47").Value = "'20"
$ws.Range("G48").Value = "'20"
$ws.Range("G49").Value = "'20"
$ws.Range("G50").Value = "'20"
$ws.Range("G51").Value = "'20"
